# Fix "Recorded By" (column G) entries: the system-generated "System"
# token was being written as the first item in the comma-separated list
# of recorders; it should instead appear last, after the real user/email
# entries (e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System").
# Rows whose "Recorded By" is just "System" alone, or "System, admin@admin.com",
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$maxRow = $usedRange.Rows.Count
$col = 7  # column G = "Recorded By"

for ($r = 1; $r -le $maxRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "

        if ($parts.Count -gt 1 -and $parts[0] -eq "System" -and $parts[1] -ne "admin@admin.com") {
            $rest = $parts[1..($parts.Count - 1)]
            $newVal = ($rest -join ", ") + ", System"
            $cell.Value = $newVal
        }
    }
}
